$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

# Slide 3 ("Fig 2a-2b") currently has, at top level:
#   1: TextBox 16
#   2: Picture 1           (the big background figure)
#   3: Group 6             (legend group, left cluster)
#   4: Group 22            (legend group, right cluster)
#
# The edit nests the two legend groups (Group 6 + Group 22) inside a new
# "Group 2", and then nests that new group together with Picture 1 inside
# another new "Group 3" - i.e. two successive Group() calls.

$legends = $s.Shapes.Range(@(3, 4))
$innerGroup = $legends.Group()

$all = $s.Shapes.Range(@(2, 3))
$outerGroup = $all.Group()
